$d = $word.ActiveDocument

# Find the "Resources:" list item that links to javatpoint - the new
# hyperlink bullet belongs to the same bulleted/numbered list, right
# after it (it's currently an empty "ListParagraph" paragraph).
$sourceIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*javatpoint*") {
        $sourceIndex = $i
    }
}

if ($sourceIndex -eq -1) {
    throw "Could not find the javatpoint list-item paragraph."
}

$sourcePara = $d.Paragraphs.Item($sourceIndex)
$targetPara = $d.Paragraphs.Item($sourceIndex + 1)

# Continue the same numbered list (numId 3) on the trailing empty
# ListParagraph paragraph instead of starting a brand new list.
$listTemplate = $sourcePara.Range.ListFormat.ListTemplate
$targetPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, 0)

# Insert the new hyperlink at the (now empty) paragraph.
$linkRange = $targetPara.Range
$linkRange.Collapse(0)
$d.Hyperlinks.Add($linkRange, "https://docs.oracle.com/javase/tutorial/jdbc/basics/index.html")

# Add the trailing space run after the hyperlink, matching the original.
$endRange = $targetPara.Range
$endRange.Collapse(0)
$endRange.InsertAfter(" ")
